$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new certificate entry (row 5)
$ws.Range("A5").Value = "Matias Sebastian Lopez Martinez"
$ws.Range("B5").Value = "Matias Sebastian Lopez Martinez"

# Update the active selection to mirror the author's last selected cell
$ws.Range("B7").Select()

$wb.Save()
